$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: B51 was stored as a text "3"; change it to a true numeric 3
$ws.Range("B51").Value = 3

# New row 52 - duplicate annotation entry appended below row 51
$ws.Range("A52").Value = "Ruilin"

# B52 must be stored as text "3" (not numeric), matching the source data
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "3"
$ws.Range("B52").Style = "Normal"

$ws.Range("C52").Value = "无"
$ws.Range("D52").Value = "FBK"
$ws.Range("E52").Value = "MET"
$ws.Range("F52").Value = "6dbc86e6-aac5-4bea-af0c-fc9177dfd16b"
$ws.Range("G52").Value = "BkJ3ibb0-_annotated.xlsx"
$ws.Range("H52").Value = "Furthermore, we have not optimized the running time of our algorithm, as it was not the focus of this work."
